$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
# Date: 11/11/2019 -> 11/19/2019 (serial 43780 -> 43788)
$ws.Range("A4").Value = 43788

# C4: bug fix text - PCF shadows now properly implemented
$ws.Range("C4").Value = "Got shadows working.  Implemented PCF shadows"

# D4: "Other Issues" - code cleanliness note (previously lived in row5/B5 slot before reshuffle)
$ws.Range("D4").Value = "Code looks messy. May cause confusion in the future. Solution: Add comments, get rid of literals, add classes, methods etc."

# G4: Supervisor comment - new content about light groupings/warnings
$ws.Range("G4").Value = "When having multiple lights make sure the are in groups of 4.  For example 1 float4 or a float3 and float etc.  Try fixing as many warnings as possible as later on it would be harder to debug if these warnings are in the way too."

# --- Row 5 ---
# Date: 11/18/2019 -> 11/26/2019 (serial 43787 -> 43795)
$ws.Range("A5").Value = 43795

# B5: new objectives text
$ws.Range("B5").Value = " Get arrays working for lights being passed over to shader.  Try fixing as many as the directx warnings as possible.  Improve program architecture further to make programming more generic.  Extra: Get assets together and setup the scene with the assets"

# Row 5 height grows to fit the longer wrapped text
$ws.Rows.Item(5).RowHeight = 135

# --- Row 6 ---
$ws.Range("A6").Value = 43802

# --- Row 7 ---
$ws.Range("A7").Value = 43809

# --- Row 8 ---
$ws.Range("A8").Value = 43816

# --- Row 9 ---
# Date cleared entirely
$ws.Range("A9").ClearContents()

# --- Sheet view: scroll so row 4 is at the top, with E4 selected ---
$ws.Range("E4").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 4
$window.ScrollColumn = 1
